$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 1.7
$ws.Range("K2").Value = 4.4
$ws.Range("N2").Value = 4.3
$ws.Range("P2").Value = 2.12
$ws.Range("Q2").Value = 1.83
$ws.Range("R2").Value = 1.43
$ws.Range("T2").Value = 1.85
$ws.Range("U2").Value = 2.06
$ws.Range("W2").Value = 2.42

# Row 3 updates
$ws.Range("N3").Value = 2.26
$ws.Range("R3").Value = 1.14
$ws.Range("S3").Value = 3
$ws.Range("T3").Value = 1.04
$ws.Range("U3").Value = 1.04

# Row 4 updates
$ws.Range("J4").Value = 4.9
$ws.Range("L4").Value = 1.39
$ws.Range("Q4").Value = 1.89
$ws.Range("S4").Value = 3.15
$ws.Range("X4").Value = 24
$ws.Range("Y4").Value = 8.800000000000001
$ws.Range("Z4").Value = 7.6
$ws.Range("AA4").Value = 11.5
$ws.Range("AB4").Value = 32
$ws.Range("AC4").Value = 11
$ws.Range("AE4").Value = 16
$ws.Range("AG4").Value = 36
$ws.Range("AH4").Value = 30
$ws.Range("AI4").Value = 980
$ws.Range("AO4").Value = 7.2
